# Update column F ("dSF") values for rows 2-36 on Sheet1
# per the commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = 2
    4  = 6
    5  = -1
    6  = 0
    7  = 1
    8  = -4
    9  = -2
    10 = 4
    11 = -3
    12 = 3
    13 = -1
    14 = -2
    15 = -1
    16 = 2
    17 = -3
    18 = 2
    19 = 3
    20 = 6
    21 = 3
    22 = -1
    23 = 4
    24 = 1
    25 = -1
    26 = -2
    27 = 7
    28 = 2
    29 = -4
    30 = 4
    31 = -1
    32 = 0
    33 = 2
    34 = -2
    35 = -4
    36 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
